# Git1.pptx -- "commit git1.pptx at 3:50pm"
#
# Target paragraph (slide 8, "Text Placeholder 1" shape) currently reads:
#   "Git " + "rm  - " + "-" + "cached filename"   (4 runs -> "Git rm  - -cached filename")
# and must become 4 separate paragraphs:
#   1) "Git rm  - -cached filename"
#   2) "" (blank)
#   3) "To show git log "
#   4) "Git log command is used"

$p = $ppt.ActivePresentation

# Locate the shape holding the "Git ... cached filename" text, searching every
# slide/shape rather than hard-coding indices, in case those ever shift.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*cached filename*") {
                $targetSlide = $sl
                $targetShape = $shp
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph whose text is the concatenation of the run-split command.
$paraIndex = -1
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    if ($tr.Paragraphs($i).Text -eq "Git rm  - -cached filename") {
        $paraIndex = $i
        break
    }
}

$para = $tr.Paragraphs($paraIndex)

# Step 1: append the new final paragraph ("Git log command is used") right after
# the target paragraph *before* touching its runs, so the new run inherits the
# formatting of the last existing run in the paragraph ("cached filename"'s rPr,
# lang="en-US" with no dirty flag) -- matching the authored diff.
$null = $para.InsertAfter("`rGit log command is used")

# Step 2: collapse the original multi-run paragraph down to a single run with
# the desired combined text. Swapping through a sentinel string that shares no
# leading/trailing characters with the original forces the engine to replace
# all runs with one run carrying the *first* original run's properties
# (lang="en-US" dirty="0"), rather than diffing prefix/suffix run-by-run.
$para = $tr.Paragraphs($paraIndex)
$para.Text = "ZZZZZZZZZZ_SENTINEL_ZZZZZZZZZZ"
$para.Text = "Git rm  - -cached filename"

# Step 3: insert the blank paragraph and the "To show git log " paragraph
# between the (now single-run) target paragraph and the final paragraph added
# in step 1. These inherit formatting from the just-collapsed run (dirty="0"),
# matching the target diff's blank-paragraph and "To show git log " paragraph.
$para = $tr.Paragraphs($paraIndex)
$null = $para.InsertAfter("`r`rTo show git log ")
